$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log-message text reused from the sheet's existing vocabulary (same
# "Info" description plus the two Message variants already present
# elsewhere in the log for these kinds of list operations).
$infoText = "Info"
$addedMsgText = "Tüm personeller listelendi"
$fullMsgText = "Kiralamalar listelendi"

# New log rows appended at the bottom of the sheet for the student-room
# placement flow: a student is placed in the room (count incremented),
# then the room becomes full, so further add attempts are logged as
# blocked (isFull = true, no more capacity).
$newRows = @(
    @{ Id = 635; Description = $infoText;  Message = $addedMsgText; Date = 45616 },
    @{ Id = 636; Description = $infoText;  Message = $fullMsgText;  Date = 45616 },
    @{ Id = 637; Description = $infoText;  Message = $fullMsgText;  Date = 45616 },
    @{ Id = 638; Description = $infoText;  Message = $fullMsgText;  Date = 45616 }
)

$startRow = 330
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value2 = $row.Id
    $ws.Cells.Item($r, 2).Value2 = $row.Description
    $ws.Cells.Item($r, 3).Value2 = $row.Message

    $dateCell = $ws.Cells.Item($r, 4)
    $dateCell.Value2 = $row.Date
    $dateCell.NumberFormat = "dd-MM-yyyy"
}
